# Guild.xlsx: unify the conception of DataNode, DataTable, Entity.
# Rename the two worksheets and make the second sheet ("DataTable") the
# active/selected tab, matching the author's re-save of the workbook.

$wb = $excel.ActiveWorkbook

$wsDataNode = $wb.Worksheets.Item(1)
$wsDataNode.Name = "DataNode"

$wsDataTable = $wb.Worksheets.Item(2)
$wsDataTable.Name = "DataTable"

# Make "DataTable" the active sheet/tab (was "Property1"/sheet1 before).
$wsDataTable.Activate()
